$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (Session 3 / pre-treatment phase measurements)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 15

# Move selection to A5, matching the post-edit state
$ws.Range("A5").Select()
